$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: Panel B (E-mini Futures) - Avg Daily Volume
$ws.Range("D26").Value = 1489883.127272727
$ws.Range("E26").Value = 814017.5469703308
$ws.Range("F26").Value = 1066294.5
$ws.Range("G26").Value = 1524449
$ws.Range("H26").Value = 1801151
$ws.Range("I26").Value = 55
$ws.Range("J26").Value = 1927055.218181818
$ws.Range("K26").Value = 651414.9672607484
$ws.Range("L26").Value = 1530795
$ws.Range("M26").Value = 1779357
$ws.Range("N26").Value = 2353539.5
$ws.Range("O26").Value = 55
$ws.Range("P26").Value = 1846542.145454546
$ws.Range("Q26").Value = 649321.2825192625
$ws.Range("R26").Value = 1398666
$ws.Range("S26").Value = 1792403
$ws.Range("T26").Value = 2294142.5
$ws.Range("U26").Value = 55
$ws.Range("V26").Value = 1677052.036363636
$ws.Range("W26").Value = 761072.4667900715
$ws.Range("X26").Value = 1382602.5
$ws.Range("Y26").Value = 1729023
$ws.Range("Z26").Value = 2221379
$ws.Range("AA26").Value = 55
$ws.Range("AB26").Value = 1639643.236363636
$ws.Range("AC26").Value = 609434.2398397655
$ws.Range("AD26").Value = 1237102
$ws.Range("AE26").Value = 1635165
$ws.Range("AF26").Value = 2005422.5
$ws.Range("AG26").Value = 55

# Row 27: Panel B (E-mini Futures) - Diff_Vol (Ann - Day)
$ws.Range("D27").Value = 356659.0181818182
$ws.Range("J27").Value = -80513.07272727272
$ws.Range("V27").Value = 169490.1090909091
$ws.Range("AB27").Value = 206898.9090909091

# Row 28: Panel B (E-mini Futures) - # Obs
$ws.Range("D28").Value = 55
$ws.Range("J28").Value = 55
$ws.Range("P28").Value = 55
$ws.Range("V28").Value = 55
$ws.Range("AB28").Value = 55
